$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 262
$ws.Range("C3").Value = 160567
$ws.Range("C4").Value = 151596
$ws.Range("C5").Value = 8971
$ws.Range("C8").Value = 64.29000000000001
